$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the _GoBack bookmark that used to sit between the "...赛龙舟"
#    run and the closing "。" run.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) Fill the existing (empty) trailing paragraph with the new date line:
#    "2022年6月7日星期二"
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$dateRange = $lastPara.Range
$dateInsertionPoint = $d.Range($dateRange.Start, $dateRange.Start)

$dateXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
    '<w:p>' + `
        '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="24"/></w:rPr><w:t>2022</w:t></w:r>' + `
        '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="24"/></w:rPr><w:t>年</w:t></w:r>' + `
        '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="24"/></w:rPr><w:t>6</w:t></w:r>' + `
        '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="24"/></w:rPr><w:t>月</w:t></w:r>' + `
        '<w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>7</w:t></w:r>' + `
        '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="24"/></w:rPr><w:t>日星期二</w:t></w:r>' + `
    '</w:p>' + `
    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$dateInsertionPoint.InsertXML($dateXml)

# ---------------------------------------------------------------------------
# 3) Append a brand-new paragraph after that date line describing the
#    weather / first day of the college entrance exam, and move the
#    _GoBack bookmark to the very end of it.
# ---------------------------------------------------------------------------
$dateParaEnd = $d.Paragraphs.Last.Range
$weatherInsertionPoint = $d.Range($dateParaEnd.End, $dateParaEnd.End)

$weatherXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
    '<w:p>' + `
        '<w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="24"/></w:rPr></w:pPr>' + `
        '<w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>晴，</w:t></w:r>' + `
        '<w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>今天是高考第一天，上午考语文，下午考数学。</w:t></w:r>' + `
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
    '</w:p>' + `
    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$weatherInsertionPoint.InsertXML($weatherXml)
